# Language.xlsx: unify the conception of DataNode, DataTable, Entity.
# The sheet that used to be called "Property1" is renamed to "DataNode".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the unified "DataNode" concept.
$ws.Name = "DataNode"

# Leave the cursor/selection where the author last left it when saving.
[void]$ws.Range("B41").Select()
